$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.121.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -7.64%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.293.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.03%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.86%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '126.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.45%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.299.77'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.83%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.466'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.24%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.32'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.45%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.115'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.41%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.367'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.39%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.869.53'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.75%  '

# Row 14
$ws.Range("E14").Value = '  -0.11%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.307.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.73%  '

# Row 16
$ws.Range("E16").Value = '  -6.84%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.69'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.79%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '59.460.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.05%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.98%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -11.24%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '346.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.77%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.547'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.99%  '

# Row 24
$ws.Range("E24").Value = '  +0.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.435.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '68.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.20%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000107'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.33%  '

# Row 28
$ws.Range("E28").Value = '  +0.08%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.00%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.63%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.69%  '

# Row 32
$ws.Range("E32").Value = '  -3.74%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.07'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.68%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.331.74'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.74%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.68%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.18'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.70'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.56%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '158.17'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.26%  '

# Row 40
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.45'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.44%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0737'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.42%  '

# Row 42
$ws.Range("E42").Value = '  +0.00%  '

# Row 43
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.98%  '

# Row 44
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.738'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.29%  '

# Row 46
$ws.Range("E46").Value = '  +3.57%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.26%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.29%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.67%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.67%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.31'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.60%  '

